$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.412.98"
$ws.Range("E2").Value = "  +4.11%  "

$ws.Range("D3").Value = "1.803.36"
$ws.Range("E3").Value = "  +1.73%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.41"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5484"
$ws.Range("E7").Value = "  +4.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3857"
$ws.Range("E8").Value = "  +6.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07619"
$ws.Range("E9").Value = "  +3.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.55"
$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.127"
$ws.Range("E11").Value = "  +3.37%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.28"
$ws.Range("E12").Value = "  +3.79%  "

$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.185"
$ws.Range("E14").Value = "  +1.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.437"
$ws.Range("E15").Value = "  +6.91%  "

$ws.Range("D16").Value = "1.805.75"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.10"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001072"
$ws.Range("E18").Value = "  +2.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06433"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.32"
$ws.Range("E21").Value = "  +3.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.973"
$ws.Range("E22").Value = "  +2.54%  "

$ws.Range("D23").Value = "28.421.31"
$ws.Range("E23").Value = "  +3.83%  "

$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.137"
$ws.Range("E25").Value = "  +2.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.96"
$ws.Range("E26").Value = "  +3.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.70"
$ws.Range("E27").Value = "  +2.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.409"
$ws.Range("E28").Value = "  +2.77%  "

$ws.Range("D29").Value = "2.011.82"
$ws.Range("E29").Value = "  +1.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.98"
$ws.Range("E30").Value = "  +2.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.124"
$ws.Range("E31").Value = "  +6.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1021"
$ws.Range("E32").Value = "  +4.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.787"
$ws.Range("E33").Value = "  +4.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.687"
$ws.Range("E34").Value = "  +1.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2320"
$ws.Range("E35").Value = "  +14.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06448"
$ws.Range("E36").Value = "  +8.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02325"
$ws.Range("E37").Value = "  +4.58%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.188"
$ws.Range("E38").Value = "  +7.25%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.817"
$ws.Range("E39").Value = "  +9.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.67"
$ws.Range("E40").Value = "  +4.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6418"
$ws.Range("E41").Value = "  +4.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.162"
$ws.Range("E42").Value = "  +1.49%  "

$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.387"
$ws.Range("E44").Value = "  -3.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.60"
$ws.Range("E45").Value = "  +4.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5986"
$ws.Range("E46").Value = "  +4.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.679"
$ws.Range("E47").Value = "  +1.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.00"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.984"
$ws.Range("E49").Value = "  +5.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.151"
$ws.Range("E50").Value = "  +3.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06892"
$ws.Range("E51").Value = "  +2.73%  "
